# Auto-generated script applying cell value updates per the commit diff.
# The workbook contains 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR);
# this script updates specific H:N cells (price/profit calculations) on each sheet
# to match the refreshed market-price data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 306804.03
$ws.Range("I15").Value = 306804.03
$ws.Range("K15").Value = 920412.0900000001
$ws.Range("M15").Value = -920243.0900000001
$ws.Range("H33").Value = 756.88464
$ws.Range("J33").Value = 1130.5
$ws.Range("L33").Value = 1130.5
$ws.Range("N33").Value = -1588.5
$ws.Range("H40").Value = 5872.8887
$ws.Range("I40").Value = 4187.9165
$ws.Range("J40").Value = 9242.833000000001
$ws.Range("K40").Value = 4187.9165
$ws.Range("L40").Value = 9242.833000000001
$ws.Range("M40").Value = -4012.9165
$ws.Range("N40").Value = -9592.833000000001
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H68").Value = 38399
$ws.Range("J68").Value = 49998
$ws.Range("L68").Value = 49998
$ws.Range("N68").Value = -51496
$ws.Range("H71").Value = 38399
$ws.Range("J71").Value = 49998
$ws.Range("L71").Value = 149994
$ws.Range("N71").Value = -157482
$ws.Range("H76").Value = 6668.6
$ws.Range("I76").Value = 5378.75
$ws.Range("K76").Value = 5378.75
$ws.Range("M76").Value = -5063.75
$ws.Range("H79").Value = 6668.6
$ws.Range("I79").Value = 5378.75
$ws.Range("K79").Value = 5378.75
$ws.Range("M79").Value = -4286.75
$ws.Range("H86").Value = 6655.8667
$ws.Range("J86").Value = 7658.5557
$ws.Range("L86").Value = 7658.5557
$ws.Range("N86").Value = -9904.555700000001
$ws.Range("H87").Value = 37064.2
$ws.Range("J87").Value = 80000
$ws.Range("L87").Value = 80000
$ws.Range("N87").Value = -82496
$ws.Range("H89").Value = 6655.8667
$ws.Range("J89").Value = 7658.5557
$ws.Range("L89").Value = 38292.7785
$ws.Range("N89").Value = -49524.7785
$ws.Range("H90").Value = 37064.2
$ws.Range("J90").Value = 80000
$ws.Range("L90").Value = 240000
$ws.Range("N90").Value = -252480
$ws.Range("H98").Value = 2153.5
$ws.Range("I98").Value = 2076.5454
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 2076.5454
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = -578.5454
$ws.Range("N98").Value = -5996
$ws.Range("H106").Value = 4259.8125
$ws.Range("I106").Value = 3705.8462
$ws.Range("K106").Value = 3705.8462
$ws.Range("M106").Value = -3074.8462
$ws.Range("H122").Value = 2153.5
$ws.Range("I122").Value = 2076.5454
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 6229.6362
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3779.6362
$ws.Range("N122").Value = -13900
$ws.Range("H124").Value = 89375
$ws.Range("J124").Value = 89375
$ws.Range("L124").Value = 89375
$ws.Range("N124").Value = -99195
$ws.Range("H125").Value = 72430.28999999999
$ws.Range("J125").Value = 126249.875
$ws.Range("L125").Value = 1136248.875
$ws.Range("N125").Value = -1141168.875
$ws.Range("H127").Value = 2016.7391
$ws.Range("I127").Value = 913.1111
$ws.Range("K127").Value = 2739.3333
$ws.Range("M127").Value = 2220.6667
$ws.Range("H131").Value = 3774
$ws.Range("I131").Value = 3500
$ws.Range("J131").Value = 3979.5
$ws.Range("K131").Value = 10500
$ws.Range("L131").Value = 11938.5
$ws.Range("M131").Value = -5460
$ws.Range("N131").Value = -22018.5
$ws.Range("H136").Value = 40779.5
$ws.Range("J136").Value = 40779.5
$ws.Range("L136").Value = 40779.5
$ws.Range("N136").Value = -50979.5
$ws.Range("H138").Value = 6156.554
$ws.Range("J138").Value = 6108.7593
$ws.Range("L138").Value = 18326.2779
$ws.Range("N138").Value = -28606.2779
$ws.Range("H141").Value = 2400.3635
$ws.Range("I141").Value = 1856
$ws.Range("J141").Value = 4850
$ws.Range("K141").Value = 5568
$ws.Range("L141").Value = 14550
$ws.Range("M141").Value = -388
$ws.Range("N141").Value = -24910

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2796.9546
$ws.Range("I2").Value = 2659.6843
$ws.Range("J2").Value = 3666.3333
$ws.Range("K2").Value = 2659.6843
$ws.Range("L2").Value = 3666.3333
$ws.Range("M2").Value = -2546.6843
$ws.Range("N2").Value = -3892.3333
$ws.Range("H32").Value = 24381
$ws.Range("I32").Value = 18563.2
$ws.Range("K32").Value = 18563.2
$ws.Range("M32").Value = -18276.2
$ws.Range("H61").Value = 11664.286
$ws.Range("I61").Value = 11056.467
$ws.Range("K61").Value = 11056.467
$ws.Range("M61").Value = -10844.467
$ws.Range("H74").Value = 4857.1787
$ws.Range("I74").Value = 3853.818
$ws.Range("J74").Value = 8536.166999999999
$ws.Range("K74").Value = 3853.818
$ws.Range("L74").Value = 8536.166999999999
$ws.Range("M74").Value = -2979.818
$ws.Range("N74").Value = -10284.167
$ws.Range("H77").Value = 4857.1787
$ws.Range("I77").Value = 3853.818
$ws.Range("J77").Value = 8536.166999999999
$ws.Range("K77").Value = 19269.09
$ws.Range("L77").Value = 42680.835
$ws.Range("M77").Value = -14901.09
$ws.Range("N77").Value = -51416.835
$ws.Range("H116").Value = 2796.9546
$ws.Range("I116").Value = 2659.6843
$ws.Range("J116").Value = 3666.3333
$ws.Range("K116").Value = 2659.6843
$ws.Range("L116").Value = 3666.3333
$ws.Range("M116").Value = -365.6842999999999
$ws.Range("N116").Value = -8254.3333
$ws.Range("H132").Value = 3545.1396
$ws.Range("I132").Value = 2840.0488
$ws.Range("K132").Value = 8520.1464
$ws.Range("M132").Value = -5990.1464
$ws.Range("H134").Value = 91304
$ws.Range("J134").Value = 91304
$ws.Range("L134").Value = 91304
$ws.Range("N134").Value = -101444
$ws.Range("H136").Value = 11664.286
$ws.Range("I136").Value = 11056.467
$ws.Range("K136").Value = 33169.401
$ws.Range("M136").Value = -30619.401

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2796.9546
$ws.Range("I3").Value = 2659.6843
$ws.Range("J3").Value = 3666.3333
$ws.Range("K3").Value = 2659.6843
$ws.Range("L3").Value = 3666.3333
$ws.Range("M3").Value = -2545.6843
$ws.Range("N3").Value = -3894.3333
$ws.Range("H86").Value = 13867.25
$ws.Range("I86").Value = 11239.833
$ws.Range("K86").Value = 11239.833
$ws.Range("M86").Value = -10116.833
$ws.Range("H89").Value = 13867.25
$ws.Range("I89").Value = 11239.833
$ws.Range("K89").Value = 56199.165
$ws.Range("M89").Value = -50583.165
$ws.Range("H105").Value = 3148.1462
$ws.Range("I105").Value = 3141.2727
$ws.Range("K105").Value = 3141.2727
$ws.Range("M105").Value = -1394.2727
$ws.Range("H107").Value = 2619.5293
$ws.Range("I107").Value = 2630.8572
$ws.Range("K107").Value = 2630.8572
$ws.Range("M107").Value = -710.8571999999999
$ws.Range("H134").Value = 5612.4043
$ws.Range("I134").Value = 4662.081
$ws.Range("K134").Value = 13986.243
$ws.Range("M134").Value = -11451.243

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8335.799999999999
$ws.Range("I31").Value = 4566.9653
$ws.Range("J31").Value = 18271.818
$ws.Range("K31").Value = 4566.9653
$ws.Range("L31").Value = 18271.818
$ws.Range("M31").Value = -4271.9653
$ws.Range("N31").Value = -18861.818
$ws.Range("H34").Value = 8335.799999999999
$ws.Range("I34").Value = 4566.9653
$ws.Range("J34").Value = 18271.818
$ws.Range("K34").Value = 4566.9653
$ws.Range("L34").Value = 18271.818
$ws.Range("M34").Value = -4364.9653
$ws.Range("N34").Value = -18675.818
$ws.Range("H64").Value = 49999
$ws.Range("J64").Value = 49999
$ws.Range("L64").Value = 49999
$ws.Range("N64").Value = -50495
$ws.Range("H67").Value = 49999
$ws.Range("J67").Value = 49999
$ws.Range("L67").Value = 49999
$ws.Range("N67").Value = -51715
$ws.Range("H69").Value = 46472.5
$ws.Range("I69").Value = 32996
$ws.Range("J69").Value = 68933.336
$ws.Range("K69").Value = 32996
$ws.Range("L69").Value = 68933.336
$ws.Range("M69").Value = -32247
$ws.Range("N69").Value = -70431.336
$ws.Range("H72").Value = 46472.5
$ws.Range("I72").Value = 32996
$ws.Range("J72").Value = 68933.336
$ws.Range("K72").Value = 98988
$ws.Range("L72").Value = 206800.008
$ws.Range("M72").Value = -95244
$ws.Range("N72").Value = -214288.008
$ws.Range("H99").Value = 6266.2915
$ws.Range("I99").Value = 7499.722
$ws.Range("J99").Value = 2566
$ws.Range("K99").Value = 7499.722
$ws.Range("L99").Value = 2566
$ws.Range("M99").Value = -6001.722
$ws.Range("N99").Value = -5562
$ws.Range("H126").Value = 6266.2915
$ws.Range("I126").Value = 7499.722
$ws.Range("J126").Value = 2566
$ws.Range("K126").Value = 22499.166
$ws.Range("L126").Value = 7698
$ws.Range("M126").Value = -20029.166
$ws.Range("N126").Value = -12638
$ws.Range("H132").Value = 3697.5938
$ws.Range("I132").Value = 3366.966
$ws.Range("K132").Value = 10100.898
$ws.Range("M132").Value = -7570.897999999999
$ws.Range("H134").Value = 3979
$ws.Range("I134").Value = 3608.1667
$ws.Range("K134").Value = 10824.5001
$ws.Range("M134").Value = -8289.500100000001
$ws.Range("H141").Value = 290868
$ws.Range("J141").Value = 347275.88
$ws.Range("L141").Value = 347275.88
$ws.Range("N141").Value = -357635.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 117.73333
$ws.Range("I2").Value = 158.6
$ws.Range("J2").Value = 36
$ws.Range("K2").Value = 951.5999999999999
$ws.Range("L2").Value = 216
$ws.Range("M2").Value = -838.5999999999999
$ws.Range("N2").Value = -442
$ws.Range("H34").Value = 5113.6665
$ws.Range("J34").Value = 10040
$ws.Range("L34").Value = 30120
$ws.Range("N34").Value = -30288
$ws.Range("H86").Value = 2461.75
$ws.Range("I86").Value = 1949.3334
$ws.Range("J86").Value = 3999
$ws.Range("K86").Value = 5848.0002
$ws.Range("L86").Value = 11997
$ws.Range("M86").Value = -4662.0002
$ws.Range("N86").Value = -14369
$ws.Range("H89").Value = 2461.75
$ws.Range("I89").Value = 1949.3334
$ws.Range("J89").Value = 3999
$ws.Range("K89").Value = 17544.0006
$ws.Range("L89").Value = 35991
$ws.Range("M89").Value = -11616.0006
$ws.Range("N89").Value = -47847
$ws.Range("H98").Value = 432.5
$ws.Range("I98").Value = 372
$ws.Range("K98").Value = 1116
$ws.Range("M98").Value = 382
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H121").Value = 724.7059
$ws.Range("I121").Value = 277.9
$ws.Range("K121").Value = 833.6999999999999
$ws.Range("M121").Value = 476.3000000000001
$ws.Range("H122").Value = 2195
$ws.Range("I122").Value = 2195
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 19755
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -17305
$ws.Range("N122").ClearContents()
$ws.Range("H130").Value = 6344.25
$ws.Range("I130").Value = 1689.5
$ws.Range("K130").Value = 5068.5
$ws.Range("M130").Value = -48.5
$ws.Range("H131").Value = 6632.857
$ws.Range("I131").Value = 2817
$ws.Range("K131").Value = 8451
$ws.Range("M131").Value = -3411
$ws.Range("H132").Value = 3731
$ws.Range("J132").Value = 3862.5
$ws.Range("L132").Value = 34762.5
$ws.Range("N132").Value = -39822.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 28999.75
$ws.Range("J44").Value = 35333
$ws.Range("L44").Value = 35333
$ws.Range("N44").Value = -36525
$ws.Range("H80").Value = 2931
$ws.Range("I80").Value = 2897
$ws.Range("J80").Value = 2948
$ws.Range("K80").Value = 2897
$ws.Range("L80").Value = 2948
$ws.Range("M80").Value = -1899
$ws.Range("N80").Value = -4944
$ws.Range("H83").Value = 2931
$ws.Range("I83").Value = 2897
$ws.Range("J83").Value = 2948
$ws.Range("K83").Value = 14485
$ws.Range("L83").Value = 14740
$ws.Range("M83").Value = -9493
$ws.Range("N83").Value = -24724
$ws.Range("H107").Value = 410.7143
$ws.Range("I107").Value = 461.8
$ws.Range("J107").Value = 283
$ws.Range("K107").Value = 461.8
$ws.Range("L107").Value = 283
$ws.Range("M107").Value = 1458.2
$ws.Range("N107").Value = -4123
$ws.Range("H113").Value = 1840.8
$ws.Range("I113").Value = 1840.8
$ws.Range("K113").Value = 1840.8
$ws.Range("M113").Value = 329.2
$ws.Range("H126").Value = 10902.607
$ws.Range("I126").Value = 10356.947
$ws.Range("K126").Value = 31070.841
$ws.Range("M126").Value = -28600.841

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 45130
$ws.Range("I35").Value = 3912.5
$ws.Range("J35").Value = 210000
$ws.Range("K35").Value = 3912.5
$ws.Range("L35").Value = 210000
$ws.Range("M35").Value = -3576.5
$ws.Range("N35").Value = -210672
$ws.Range("H40").Value = 7673.25
$ws.Range("I40").Value = 5564.3335
$ws.Range("J40").Value = 14000
$ws.Range("K40").Value = 5564.3335
$ws.Range("L40").Value = 14000
$ws.Range("M40").Value = -5428.3335
$ws.Range("N40").Value = -14272
$ws.Range("H46").Value = 4226.1177
$ws.Range("I46").Value = 4493.875
$ws.Range("J46").Value = 3988.111
$ws.Range("K46").Value = 4493.875
$ws.Range("L46").Value = 3988.111
$ws.Range("M46").Value = -4305.875
$ws.Range("N46").Value = -4364.111
$ws.Range("H55").Value = 670.17645
$ws.Range("I55").Value = 310.3
$ws.Range("J55").Value = 1184.2858
$ws.Range("K55").Value = 310.3
$ws.Range("L55").Value = 1184.2858
$ws.Range("M55").Value = -137.3
$ws.Range("N55").Value = -1530.2858
$ws.Range("H122").Value = 9961.25
$ws.Range("I122").Value = 5845
$ws.Range("J122").Value = 11333.333
$ws.Range("K122").Value = 17535
$ws.Range("L122").Value = 33999.999
$ws.Range("M122").Value = -15085
$ws.Range("N122").Value = -38899.999
$ws.Range("H124").Value = 75389.39999999999
$ws.Range("J124").Value = 75389.39999999999
$ws.Range("L124").Value = 75389.39999999999
$ws.Range("N124").Value = -85209.39999999999
$ws.Range("H132").Value = 3753.0645
$ws.Range("I132").Value = 3753.0645
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11259.1935
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8729.193499999999
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 9440.75
$ws.Range("I136").Value = 6469
$ws.Range("K136").Value = 19407
$ws.Range("M136").Value = -16857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 40831.332
$ws.Range("I54").Value = 49999
$ws.Range("J54").Value = 36247.5
$ws.Range("K54").Value = 49999
$ws.Range("L54").Value = 36247.5
$ws.Range("M54").Value = -49479
$ws.Range("N54").Value = -37287.5
$ws.Range("H62").Value = 14000
$ws.Range("J62").Value = 14000
$ws.Range("L62").Value = 14000
$ws.Range("N62").Value = -15248
$ws.Range("H65").Value = 14000
$ws.Range("J65").Value = 14000
$ws.Range("L65").Value = 70000
$ws.Range("N65").Value = -76240
$ws.Range("H75").Value = 98588.28999999999
$ws.Range("I75").Value = 95059
$ws.Range("K75").Value = 95059
$ws.Range("M75").Value = -94123
$ws.Range("H78").Value = 98588.28999999999
$ws.Range("I78").Value = 95059
$ws.Range("K78").Value = 285177
$ws.Range("M78").Value = -280497
$ws.Range("H126").Value = 2852.9546
$ws.Range("I126").Value = 2107.875
$ws.Range("K126").Value = 6323.625
$ws.Range("M126").Value = -3853.625
$ws.Range("H132").Value = 7606.148
$ws.Range("I132").Value = 7214.68
$ws.Range("J132").Value = 12499.5
$ws.Range("K132").Value = 21644.04
$ws.Range("L132").Value = 37498.5
$ws.Range("M132").Value = -19114.04
$ws.Range("N132").Value = -42558.5
$ws.Range("H136").Value = 7488.0586
$ws.Range("I136").Value = 6524.3335
$ws.Range("J136").Value = 9801
$ws.Range("K136").Value = 19573.0005
$ws.Range("L136").Value = 29403
$ws.Range("M136").Value = -17023.0005
$ws.Range("N136").Value = -34503
